$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.792.97'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.160.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.31%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.97'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.36'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.11%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.157.94'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.62%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.507'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +13.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.05'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.671.02'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.883.24'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.26'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +7.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.162.66'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.73%  '

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '520.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.54%  '

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.111'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.739'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.86'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.38'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.00%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +9.61%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.07'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.28%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.55%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +10.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.60'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.79'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '489.11'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0868'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0424'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.01'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.117.41'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.68'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.299'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +14.46%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +15.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.31'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0582'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +13.02%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.50'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.14%  '
